# Update the lattice-multiplication practice table: replace every cell's
# multiplication problem (all 15 cells across the 5x3 table) with the new
# problem set, preserving the "digits / ---- / lattice box" layout (each
# line inside a cell is a separate <w:t> run joined by <w:br/>, represented
# here by the vertical-tab character so Word keeps them as line breaks).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11   # vertical tab -> becomes <w:br/> between the w:t runs

function Set-LatticeCell([int]$row, [int]$col, [string[]]$lines) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = [string]::Join($vt, $lines)
}

# Row 1, Col 1: 22 x 65 -> 19 x 94
Set-LatticeCell 1 1 @("19 x 94", "  9    4", "  ----", "1|    |", "9|    |")

# Row 1, Col 2: 46 x 69 -> 78 x 58
Set-LatticeCell 1 2 @("78 x 58", "  5    8", "  ----", "7|    |", "8|    |")

# Row 1, Col 3: 10 x 26 -> 74 x 12
Set-LatticeCell 1 3 @("74 x 12", "  1    2", "  ----", "7|    |", "4|    |")

# Row 2, Col 1: 59 x 16 -> 46 x 27
Set-LatticeCell 2 1 @("46 x 27", "  2    7", "  ----", "4|    |", "6|    |")

# Row 2, Col 2: 41 x 11 -> 37 x 33
Set-LatticeCell 2 2 @("37 x 33", "  3    3", "  ----", "3|    |", "7|    |")

# Row 2, Col 3: 51 x 98 -> 79 x 26
Set-LatticeCell 2 3 @("79 x 26", "  2    6", "  ----", "7|    |", "9|    |")

# Row 3, Col 1: 68 x 57 -> 55 x 99
Set-LatticeCell 3 1 @("55 x 99", "  9    9", "  ----", "5|    |", "5|    |")

# Row 3, Col 2: 74 x 62 -> 72 x 26
Set-LatticeCell 3 2 @("72 x 26", "  2    6", "  ----", "7|    |", "2|    |")

# Row 3, Col 3: 46 x 97 -> 44 x 60
Set-LatticeCell 3 3 @("44 x 60", "  6    0", "  ----", "4|    |", "4|    |")

# Row 4, Col 1: 27 x 71 -> 83 x 46
Set-LatticeCell 4 1 @("83 x 46", "  4    6", "  ----", "8|    |", "3|    |")

# Row 4, Col 2: 10 x 94 -> 14 x 61
Set-LatticeCell 4 2 @("14 x 61", "  6    1", "  ----", "1|    |", "4|    |")

# Row 4, Col 3: 79 x 26 -> 50 x 20
Set-LatticeCell 4 3 @("50 x 20", "  2    0", "  ----", "5|    |", "0|    |")

# Row 5, Col 1: 53 x 37 -> 96 x 45
Set-LatticeCell 5 1 @("96 x 45", "  4    5", "  ----", "9|    |", "6|    |")

# Row 5, Col 2: 72 x 62 -> 22 x 62
Set-LatticeCell 5 2 @("22 x 62", "  6    2", "  ----", "2|    |", "2|    |")

# Row 5, Col 3: 99 x 64 -> 86 x 11
Set-LatticeCell 5 3 @("86 x 11", "  1    1", "  ----", "8|    |", "6|    |")
